$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Rushing" - Week 16 stat updates for existing rows
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: J.Hurts
$rushing.Range("C2").Value = 44
$rushing.Range("F2").Value = 30

# Row 4: M.Sanders
$rushing.Range("C4").Value = 46
$rushing.Range("D4").Value = 24
$rushing.Range("F4").Value = 11

# Row 5: B.Scott
$rushing.Range("C5").Value = 50
$rushing.Range("D5").Value = 27
$rushing.Range("E5").Value = 13
$rushing.Range("F5").Value = 17

# Row 7: J.Howard
$rushing.Range("C7").Value = 37
$rushing.Range("D7").Value = 29
$rushing.Range("F7").Value = 21

# ---------------------------------------------------------------------
# Sheet "Receiving" - Week 16 stat updates + new row for J.Howard (week
# 16 only) + a new trailing row for T.Jackson (season sim from Wk 17)
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Shift the B-column (shared-string) player labels down by one row,
# working bottom-up so no label is ever left with zero references.
$receiving.Range("B13").Value = "T.Jackson"
$receiving.Range("B12").Value = "J.Stoll"
$receiving.Range("B11").Value = "D.Goedert"
$receiving.Range("B10").Value = "J.Arcega-Whiteside"
$receiving.Range("B9").Value = "G.Ward"
$receiving.Range("B8").Value = "Q.Watkins"
$receiving.Range("B7").Value = "J.Reagor"
$receiving.Range("B6").Value = "D.Smith"
$receiving.Range("B5").Value = "J.Howard"

# Give new row 13 the same formatting as row 12 (bold / centered /
# top-aligned / thin border matching the rest of column A).
$receiving.Range("A12").Copy()
$receiving.Range("A13").PasteSpecial(-4122)

# Column A index (0-based running counter)
$receiving.Range("A5").Value = 3
$receiving.Range("A6").Value = 4
$receiving.Range("A7").Value = 5
$receiving.Range("A8").Value = 6
$receiving.Range("A9").Value = 7
$receiving.Range("A10").Value = 8
$receiving.Range("A11").Value = 9
$receiving.Range("A12").Value = 10
$receiving.Range("A13").Value = 11

# Row 5: J.Howard (new entry, week 16 stats only)
$receiving.Range("C5").Value = 4
$receiving.Range("D5").Value = 2
$receiving.Range("E5").Value = 0
$receiving.Range("F5").Value = 0
$receiving.Range("G5").Value = 1
$receiving.Range("H5").Value = 1

# Row 6: D.Smith
$receiving.Range("C6").Value = 62
$receiving.Range("D6").Value = 45
$receiving.Range("E6").Value = 32
$receiving.Range("F6").Value = 14
$receiving.Range("G6").Value = 8
$receiving.Range("H6").Value = 5

# Row 7: J.Reagor
$receiving.Range("C7").Value = 41
$receiving.Range("D7").Value = 27
$receiving.Range("E7").Value = 12
$receiving.Range("F7").Value = 4
$receiving.Range("G7").Value = 5
$receiving.Range("H7").Value = 3

# Row 8: Q.Watkins
$receiving.Range("C8").Value = 32
$receiving.Range("D8").Value = 25
$receiving.Range("E8").Value = 18
$receiving.Range("F8").Value = 10
$receiving.Range("G8").Value = 9
$receiving.Range("H8").Value = 5

# Row 9: G.Ward
$receiving.Range("C9").Value = 8
$receiving.Range("D9").Value = 4
$receiving.Range("E9").Value = 0
$receiving.Range("F9").Value = 0
$receiving.Range("G9").Value = 7
$receiving.Range("H9").Value = 3

# Row 10: J.Arcega-Whiteside
$receiving.Range("C10").Value = 1
$receiving.Range("D10").Value = 1
$receiving.Range("E10").Value = 1
$receiving.Range("F10").Value = 1
$receiving.Range("G10").Value = 0
$receiving.Range("H10").Value = 0

# Row 11: D.Goedert
$receiving.Range("C11").Value = 50
$receiving.Range("D11").Value = 34
$receiving.Range("E11").Value = 18
$receiving.Range("F11").Value = 15
$receiving.Range("G11").Value = 7
$receiving.Range("H11").Value = 5

# Row 12: J.Stoll
$receiving.Range("C12").Value = 4
$receiving.Range("D12").Value = 3
$receiving.Range("E12").Value = 0
$receiving.Range("F12").Value = 0
$receiving.Range("G12").Value = 1
$receiving.Range("H12").Value = 1

# Row 13: T.Jackson (new trailing row)
$receiving.Range("C13").Value = 1
$receiving.Range("D13").Value = 0
$receiving.Range("E13").Value = 1
$receiving.Range("F13").Value = 0
$receiving.Range("G13").Value = 0
$receiving.Range("H13").Value = 0
